# Add "My profile" page entries (row 42-44) to the "Journal De Bord" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal De Bord")

# --- New row 42: date (A), time (B), description (C) ---
# Copy formats from existing styled cells so no new cellXfs are created,
# then overwrite with the real values/text.
$ws.Range("A32").Copy($ws.Range("A42"))
$ws.Range("A42").Value = 43560

$ws.Range("B32").Copy($ws.Range("B42"))
$ws.Range("B42").Value = 0.3125

$ws.Range("C42").Value = "Analyse du critère 2"

# --- New row 43: time (B), description (C) ---
$ws.Range("B33").Copy($ws.Range("B43"))
$ws.Range("B43").Value = 0.34722222222222227

$ws.Range("C43").Value = "Implémentation des information personnel"

# --- New row 44: an empty, time-styled cell (matches the trailing blank row) ---
$ws.Range("B26").Copy($ws.Range("B44"))

# --- Column width tweaks (widened to fit the new, longer text) ---
$ws.Columns.Item(2).ColumnWidth = 6.7109375
$ws.Columns.Item(3).ColumnWidth = 42.140625

# --- View: scroll down and move the selection to the newly added cell ---
$ws.Activate()
$ws.Range("C43").Select()
